$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.202.66"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.99%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.819.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.08%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.013"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +1.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -8.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.014"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5089"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2390"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -25.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.05853"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -13.78%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.842.57"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06791"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -12.32%  "
$ws.Range("B12").Value = "Litecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "78.70"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -10.28%  "
$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.06"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -25.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.369"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -12.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5679"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -27.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.018"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.67%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.014"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.229.57"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.98%  "
$ws.Range("B19").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C19").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.057.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.67%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -23.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000005925"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -25.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.850"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -16.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.048"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -15.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.741"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -17.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "129.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.432"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -15.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.795"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -17.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "13.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -18.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "96.65"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -13.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08232"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.545"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -14.70%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.770"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04169"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -14.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.081"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -24.35%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.032"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -9.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.976"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.76%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6024"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -17.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.046"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -8.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.015"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.43%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8267"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.45%  "
$ws.Range("B41").Value = "Quant"
$ws.Range("C41").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.48"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -10.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.365"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -9.47%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.01414"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -19.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3576"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -25.12%  "
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.05262"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -10.07%  "
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.940"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -22.56%  "
$ws.Range("B47").Value = "Elrond"
$ws.Range("C47").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "29.82"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -14.41%  "
$ws.Range("B48").Value = "USDD"
$ws.Range("C48").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.010"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "51.97"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -12.70%  "
$ws.Range("B50").Value = "TrueUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.011"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.10%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.299"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -18.78%  "
